# Update "想去人数" (F column) counts to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1)
$ws1.Range("F3").Value  = 545
$ws1.Range("F8").Value  = 45
$ws1.Range("F11").Value = 1179
$ws1.Range("F14").Value = 835
$ws1.Range("F15").Value = 847
$ws1.Range("F18").Value = 72
$ws1.Range("F20").Value = 725
$ws1.Range("F21").Value = 1726
$ws1.Range("F22").Value = 2609
$ws1.Range("F23").Value = 740
$ws1.Range("F25").Value = 2012
$ws1.Range("F26").Value = 472
$ws1.Range("F27").Value = 2897
$ws1.Range("F28").Value = 539
$ws1.Range("F31").Value = 138
$ws1.Range("F34").Value = 1014
$ws1.Range("F35").Value = 1732
$ws1.Range("F38").Value = 545
$ws1.Range("F39").Value = 171
$ws1.Range("F40").Value = 126
$ws1.Range("F42").Value = 28

# 演出 sheet (sheet2)
$ws2.Range("F3").Value  = 138
$ws2.Range("F8").Value  = 3
$ws2.Range("F12").Value = 75

# 全部类型 sheet (sheet4)
$ws4.Range("F4").Value  = 545
$ws4.Range("F9").Value  = 45
$ws4.Range("F12").Value = 1179
$ws4.Range("F14").Value = 835
$ws4.Range("F15").Value = 847
$ws4.Range("F16").Value = 138
$ws4.Range("F20").Value = 72
$ws4.Range("F21").Value = 725
$ws4.Range("F22").Value = 1726
$ws4.Range("F23").Value = 2609
$ws4.Range("F24").Value = 740
$ws4.Range("F28").Value = 2897
$ws4.Range("F29").Value = 539
$ws4.Range("F30").Value = 3
$ws4.Range("F35").Value = 75
$ws4.Range("F37").Value = 138
$ws4.Range("F40").Value = 1014
$ws4.Range("F41").Value = 1732
$ws4.Range("F44").Value = 545
$ws4.Range("F45").Value = 171
$ws4.Range("F46").Value = 126
$ws4.Range("F48").Value = 28
